# Add victim age and offense category/groupings for demos
#
# Adds a new "OffenseCategory1" column (D) to the UCROffenseCodeType sheet,
# classifying each UCR offense code as "Serious Violence", "Simple Assault"
# or "Other", then makes that sheet the active tab/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UCROffenseCodeType")
$prev = $wb.Worksheets.Item("ResidentStatusOfPersonType")

# Column header
$ws.Range("D1").Value = "OffenseCategory1"

# Row -> category, keyed by the UCROffenseCodeTypeID in column A.
# Note: "Simple Assault" (row 8) is written before "Serious Violence"
# (row 2 et al.) so that new shared strings are interned in that order.
$ws.Cells.Item(8, 4).Value = "Simple Assault"    # 132 / 13B - Assault-Simple Assault

$ws.Cells.Item(2, 4).Value = "Serious Violence"  # 91  / 09A - Homicide-Murder and Nonnegligent Manslaughter
$ws.Cells.Item(5, 4).Value = "Serious Violence"  # 100 / 100 - Kidnapping/Abduction
$ws.Cells.Item(6, 4).Value = "Serious Violence"  # 120 / 120 - Robbery
$ws.Cells.Item(7, 4).Value = "Serious Violence"  # 131 / 13A - Assault-Aggravated Assault
$ws.Cells.Item(40, 4).Value = "Serious Violence" # 111 / 11A - Sex Offenses, Forcible-Forcible Rape

$otherRows = 3,4,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,41,42,43,44,45,46,47,48,49,50,51,52,53,54,55,56,57,58
foreach ($r in $otherRows) {
    $ws.Cells.Item($r, 4).Value = "Other"
}

# New column width
$ws.Columns.Item(4).ColumnWidth = 16.75

# UCROffenseCodeType becomes the active sheet/tab, with the given selection;
# ResidentStatusOfPersonType (previously active) loses its selection mark.
$prev.Range("A4").Select()
$ws.Activate()
$ws.Range("E45").Select()
